$d = $word.ActiveDocument

# --- Replace the text of paragraphs 1-9 (1-indexed) while keeping each paragraph mark intact ---
$p0 = $d.Paragraphs.Item(1)
$r0 = $p0.Range
$r0.MoveEnd(1, -1) | Out-Null
$r0.Text = 'המאמר היומי של מייק: 29.08.25'

$p1 = $d.Paragraphs.Item(2)
$r1 = $p1.Range
$r1.MoveEnd(1, -1) | Out-Null
$r1.Text = 'DEEP THINK WITH CONFIDENCE'

$p2 = $d.Paragraphs.Item(3)
$r2 = $p2.Range
$r2.MoveEnd(1, -1) | Out-Null
$r2.Text = 'סקירה מס 497 - עוד 3 סקירת בדרך ל-500 והיום סקירה קצרה של מאמר בעל שם מפוצץ (שאכן נהנה מהייפ משמעותי) עם רעיון די אינטואיטיבי שגרם לי לתהות איך אף אחד לא עשה את זה קודם (אם זה נכון). המאמר מציע שיטה מבוססת אנטרופיה לדגימה ממודלי שפה אוטורגרסיביים (למרות שלדעתי אפשר יחסית בקלות להרחיב את הגישה המוצעת למודלים שמגנרטים פלט בצורה לא אוטורגרסיבית כמו מודלי שפה מבוססי דיפוזיה). כמו שאתם בטח יודעים אנטרופיה הינו מדד לאי ודאות וניתן להשתמש בו במודלי שפה למטרת שערוך של ״מידת הביטחון״ של המודל בפלט שהוא מג''נרט.'

$p3 = $d.Paragraphs.Item(4)
$r3 = $p3.Range
$r3.MoveEnd(1, -1) | Out-Null
$r3.Text = 'מודלי שפה אוטורגרסיביים מגנרטים כל טוקן בהתבסס על ההתפלגות של הטוקן הזה בהינתן ההקשר הקודם לו. ככל שהאנטרופיה של הטוקן הנחזה, השווה למינוס לוג של ההסתברות שלו, גבוהה יותר האי ודאות שלו גבוה יותר. כלומר ככל שההסתברות הטוקן יורדת, אי הוודאות עולה הקשורה בבחירותו עולה. כאמור המחברים מציעים שיטת דגימה מבוססת אנטרופיה ממוצעת של הטוקנים בטקסט מג''ונרט.'

$p4 = $d.Paragraphs.Item(5)
$r4 = $p4.Range
$r4.MoveEnd(1, -1) | Out-Null
$r4.Text = 'בפרט אם במקרים שהמודל מגנרט כמה כמה תשובות לשאלה מתמטית ואז אנו בוחרים את התשובה הנכונה לא עם ה-majority vote פשוט (כלומר התשובה הסופית שרוב התשובות התכנסו אליה) אלא על ידי משקול תשובה עם הודאות שלה כלומר עם הממוצע של האנטרופיה של כל הטוקנים שלה. כך תשובות שהמודל ממש לא בטוח בהם מפולטרות. המחברים גם מציעים לקבוע סף של אי וודאות מקסימלית של תשובת המודל. אם אי הוודאות הממוצעת השוטפת של התשובה (מחושבת מחדש עבור כל טוקן מג''ונרט), התשובה נפסלת והמודל מספיק לגנרט אותה. הסף נקבע בתור אחוזון של אי הוודאויות של התשובות הנכונות בשלב ה-warmup.'

$p5 = $d.Paragraphs.Item(6)
$r5 = $p5.Range
$r5.MoveEnd(1, -1) | Out-Null
$r5.Text = 'בנוסף המחברים מציעים לקבוע את מספר תשובות הנדגמות מהמודל בהתבסס על קושי השאלה. ככל ש״יש מעט מדי הסכמה״ בין התוצאות של התשובות השונות המודל מגנרט יותר תשובות כאשר התשובות בעלות אי וודאות גבוהה מדי מפולטרות כאמור.'

$p6 = $d.Paragraphs.Item(7)
$r6 = $p6.Range
$r6.MoveEnd(1, -1) | Out-Null
$r6.Text = 'מאמר נחמד אבל משאיר תחושה שכבר ראיתי משהו כזה בעבר…. '

$p7 = $d.Paragraphs.Item(8)
$r7 = $p7.Range
$r7.MoveEnd(1, -1) | Out-Null
$r7.Text = 'https://arxiv.org/abs/2508.15260 '

$p8 = $d.Paragraphs.Item(9)
$r8 = $p8.Range
$r8.MoveEnd(1, -1) | Out-Null
$r8.Text = 'שחר ה-AI המקיאווליאני? חשיפת התבונה האסטרטגית של LLMS'

# --- Remove paragraphs 10 through 20 (old indices 9..19), which no longer exist in the new version ---
$firstToRemove = $d.Paragraphs.Item(10)
$lastToRemove = $d.Paragraphs.Item($d.Paragraphs.Count)
$rangeToDelete = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
$rangeToDelete.Delete()

Write-Host 'Final paragraph count:' $d.Paragraphs.Count